$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4375153333333333
$ws.Range("H2").Value = 1.312546
$ws.Range("I2").Value = 0.002535486401940996
$ws.Range("J2").Value = 0.002555908833496712
$ws.Range("M2").Value = 0.600843
$ws.Range("N2").Value = 1.802529
$ws.Range("O2").Value = 0.1216566842860732
$ws.Range("P2").Value = 0.148762828988552
$ws.Range("Q2").Value = 0.262878025426
$ws.Range("R2").Value = 2.365902228834
$ws.Range("S2").Value = 0.0003084588687125674
$ws.Range("T2").Value = 0.0003802242287078008
$ws.Range("G3").Value = 0.4375153333333333
$ws.Range("H3").Value = 1.312546
$ws.Range("I3").Value = 0.002535486401940996
$ws.Range("J3").Value = 0.002555908833496712
$ws.Range("O3").Value = 0.2991130341144489
$ws.Range("P3").Value = 0.3657579639239645
$ws.Range("Q3").Value = 0.6463290056655555
$ws.Range("R3").Value = 5.816961050990001
$ws.Range("S3").Value = 0.0007583970306404984
$ws.Range("T3").Value = 0.0009348440109150324
$ws.Range("G4").Value = 0.4375153333333333
$ws.Range("H4").Value = 1.312546
$ws.Range("I4").Value = 0.002535486401940996
$ws.Range("J4").Value = 0.002555908833496712
$ws.Range("M4").Value = 0.111967
$ws.Range("N4").Value = 0.335901
$ws.Range("O4").Value = 0.02267070427625646
$ws.Range("P4").Value = 0.02772193014375004
$ws.Range("Q4").Value = 0.04898727932733333
$ws.Range("R4").Value = 0.440885513946
$ws.Range("S4").Value = 0.00005748126241487383
$ws.Range("T4").Value = 0.0000708547261359895
$ws.Range("G5").Value = 0.4375153333333333
$ws.Range("H5").Value = 1.312546
$ws.Range("I5").Value = 0.002535486401940996
$ws.Range("J5").Value = 0.002555908833496712
$ws.Range("M5").Value = 2.6997255
$ws.Range("N5").Value = 5.399451
$ws.Range("O5").Value = 0.5466314042313235
$ws.Range("P5").Value = 0.4456170223863617
$ws.Range("Q5").Value = 1.181171302041
$ws.Range("R5").Value = 7.087027812246
$ws.Range("S5").Value = 0.001385976492302433
$ws.Range("T5").Value = 0.001138956483873804
$ws.Range("G6").Value = 0.4375153333333333
$ws.Range("H6").Value = 1.312546
$ws.Range("I6").Value = 0.002535486401940996
$ws.Range("J6").Value = 0.002555908833496712
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.04903366666666667
$ws.Range("N6").Value = 0.147101
$ws.Range("O6").Value = 0.009928173091897913
$ws.Range("P6").Value = 0.01214025455737189
$ws.Range("Q6").Value = 0.02145298101622222
$ws.Range("R6").Value = 0.193076829146
$ws.Range("S6").Value = 0.00002517274787062365
$ws.Range("T6").Value = 0.00003102938386408552
$ws.Range("I7").Value = 0.005820447907772805
$ws.Range("J7").Value = 0.005867329523437988
$ws.Range("M7").Value = 0.600843
$ws.Range("N7").Value = 1.802529
$ws.Range("O7").Value = 0.1216566842860732
$ws.Range("P7").Value = 0.148762828988552
$ws.Range("Q7").Value = 0.6034612735130001
$ws.Range("R7").Value = 5.431151461617
$ws.Range("S7").Value = 0.0007080963935194515
$ws.Range("T7").Value = 0.000872840538514688
$ws.Range("I8").Value = 0.005820447907772805
$ws.Range("J8").Value = 0.005867329523437988
$ws.Range("O8").Value = 0.2991130341144489
$ws.Range("P8").Value = 0.3657579639239645
$ws.Range("S8").Value = 0.00174097183359902
$ws.Range("T8").Value = 0.002146022500163644
$ws.Range("I9").Value = 0.005820447907772805
$ws.Range("J9").Value = 0.005867329523437988
$ws.Range("M9").Value = 0.111967
$ws.Range("N9").Value = 0.335901
$ws.Range("O9").Value = 0.02267070427625646
$ws.Range("P9").Value = 0.02772193014375004
$ws.Range("Q9").Value = 0.1124549148636667
$ws.Range("R9").Value = 1.012094233773
$ws.Range("S9").Value = 0.0001319536532724729
$ws.Range("T9").Value = 0.0001626536991791101
$ws.Range("I10").Value = 0.005820447907772805
$ws.Range("J10").Value = 0.005867329523437988
$ws.Range("M10").Value = 2.6997255
$ws.Range("N10").Value = 5.399451
$ws.Range("O10").Value = 0.5466314042313235
$ws.Range("P10").Value = 0.4456170223863617
$ws.Range("Q10").Value = 2.7114900038205
$ws.Range("R10").Value = 16.268940022923
$ws.Range("S10").Value = 0.003181639613081118
$ws.Range("T10").Value = 0.002614581911594027
$ws.Range("I11").Value = 0.005820447907772805
$ws.Range("J11").Value = 0.005867329523437988
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.04903366666666667
$ws.Range("N11").Value = 0.147101
$ws.Range("O11").Value = 0.009928173091897913
$ws.Range("P11").Value = 0.01214025455737189
$ws.Range("Q11").Value = 0.04924733904144445
$ws.Range("R11").Value = 0.443226051373
$ws.Range("S11").Value = 0.00005778641430074347
$ws.Range("T11").Value = 0.00007123087398652068
$ws.Range("G12").Value = 99.58055866666666
$ws.Range("H12").Value = 298.741676
$ws.Range("I12").Value = 0.577088694179909
$ws.Range("J12").Value = 0.5817369361698658
$ws.Range("M12").Value = 0.600843
$ws.Range("N12").Value = 1.802529
$ws.Range("O12").Value = 0.1216566842860732
$ws.Range("P12").Value = 0.148762828988552
$ws.Range("Q12").Value = 59.832281610956
$ws.Range("R12").Value = 538.490534498604
$ws.Range("S12").Value = 0.07020669707290744
$ws.Range("T12").Value = 0.08654083235176195
$ws.Range("G13").Value = 99.58055866666666
$ws.Range("H13").Value = 298.741676
$ws.Range("I13").Value = 0.577088694179909
$ws.Range("J13").Value = 0.5817369361698658
$ws.Range("O13").Value = 0.2991130341144489
$ws.Range("P13").Value = 0.3657579639239645
$ws.Range("Q13").Value = 147.1075378691044
$ws.Range("R13").Value = 1323.96784082194
$ws.Range("S13").Value = 0.1726147502692979
$ws.Range("T13").Value = 0.2127749173128554
$ws.Range("G14").Value = 99.58055866666666
$ws.Range("H14").Value = 298.741676
$ws.Range("I14").Value = 0.577088694179909
$ws.Range("J14").Value = 0.5817369361698658
$ws.Range("M14").Value = 0.111967
$ws.Range("N14").Value = 0.335901
$ws.Range("O14").Value = 0.02267070427625646
$ws.Range("P14").Value = 0.02772193014375004
$ws.Range("Q14").Value = 11.14973641223067
$ws.Range("R14").Value = 100.347627710076
$ws.Range("S14").Value = 0.01308300712692372
$ws.Range("T14").Value = 0.01612687070654019
$ws.Range("G15").Value = 99.58055866666666
$ws.Range("H15").Value = 298.741676
$ws.Range("I15").Value = 0.577088694179909
$ws.Range("J15").Value = 0.5817369361698658
$ws.Range("M15").Value = 2.6997255
$ws.Range("N15").Value = 5.399451
$ws.Range("O15").Value = 0.5466314042313235
$ws.Range("P15").Value = 0.4456170223863617
$ws.Range("Q15").Value = 268.840173536646
$ws.Range("R15").Value = 1613.041041219876
$ws.Range("S15").Value = 0.3154548032655845
$ws.Range("T15").Value = 0.2592318813081806
$ws.Range("G16").Value = 99.58055866666666
$ws.Range("H16").Value = 298.741676
$ws.Range("I16").Value = 0.577088694179909
$ws.Range("J16").Value = 0.5817369361698658
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.04903366666666667
$ws.Range("N16").Value = 0.147101
$ws.Range("O16").Value = 0.009928173091897913
$ws.Range("P16").Value = 0.01214025455737189
$ws.Range("Q16").Value = 4.882799920141778
$ws.Range("R16").Value = 43.945199281276
$ws.Range("S16").Value = 0.005729436445195477
$ws.Range("T16").Value = 0.007062434490527773
$ws.Range("G17").Value = 4.1363315
$ws.Range("H17").Value = 8.272663
$ws.Range("I17").Value = 0.02397084507248554
$ws.Range("J17").Value = 0.01610928107528529
$ws.Range("M17").Value = 0.600843
$ws.Range("N17").Value = 1.802529
$ws.Range("O17").Value = 0.1216566842860732
$ws.Range("P17").Value = 0.148762828988552
$ws.Range("Q17").Value = 2.4852858274545
$ws.Range("R17").Value = 14.911714964727
$ws.Range("S17").Value = 0.002916213531053747
$ws.Range("T17").Value = 0.002396462225731184
$ws.Range("G18").Value = 4.1363315
$ws.Range("H18").Value = 8.272663
$ws.Range("I18").Value = 0.02397084507248554
$ws.Range("J18").Value = 0.01610928107528529
$ws.Range("O18").Value = 0.2991130341144489
$ws.Range("P18").Value = 0.3657579639239645
$ws.Range("Q18").Value = 6.110485328890833
$ws.Range("R18").Value = 36.662911973345
$ws.Range("S18").Value = 0.007169992199918537
$ws.Range("T18").Value = 0.005892097846375201
$ws.Range("G19").Value = 4.1363315
$ws.Range("H19").Value = 8.272663
$ws.Range("I19").Value = 0.02397084507248554
$ws.Range("J19").Value = 0.01610928107528529
$ws.Range("M19").Value = 0.111967
$ws.Range("N19").Value = 0.335901
$ws.Range("O19").Value = 0.02267070427625646
$ws.Range("P19").Value = 0.02772193014375004
$ws.Range("Q19").Value = 0.4631326290605
$ws.Range("R19").Value = 2.778795774363
$ws.Range("S19").Value = 0.000543435939890279
$ws.Range("T19").Value = 0.0004465803646350934
$ws.Range("G20").Value = 4.1363315
$ws.Range("H20").Value = 8.272663
$ws.Range("I20").Value = 0.02397084507248554
$ws.Range("J20").Value = 0.01610928107528529
$ws.Range("M20").Value = 2.6997255
$ws.Range("N20").Value = 5.399451
$ws.Range("O20").Value = 0.5466314042313235
$ws.Range("P20").Value = 0.4456170223863617
$ws.Range("Q20").Value = 11.16695962700325
$ws.Range("R20").Value = 44.667838508013
$ws.Range("S20").Value = 0.01310321670258427
$ws.Range("T20").Value = 0.007178569865553599
$ws.Range("G21").Value = 4.1363315
$ws.Range("H21").Value = 8.272663
$ws.Range("I21").Value = 0.02397084507248554
$ws.Range("J21").Value = 0.01610928107528529
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.04903366666666667
$ws.Range("N21").Value = 0.147101
$ws.Range("O21").Value = 0.009928173091897913
$ws.Range("P21").Value = 0.01214025455737189
$ws.Range("Q21").Value = 0.2028194999938333
$ws.Range("R21").Value = 1.216916999963
$ws.Range("S21").Value = 0.0002379866990387046
$ws.Range("T21").Value = 0.000195570772990217
$ws.Range("G22").Value = 67.39800266666667
$ws.Range("H22").Value = 202.194008
$ws.Range("I22").Value = 0.3905845264378918
$ws.Range("J22").Value = 0.3937305443979143
$ws.Range("M22").Value = 0.600843
$ws.Range("N22").Value = 1.802529
$ws.Range("O22").Value = 0.1216566842860732
$ws.Range("P22").Value = 0.148762828988552
$ws.Range("Q22").Value = 40.495618116248
$ws.Range("R22").Value = 364.460563046232
$ws.Range("S22").Value = 0.04751721841988001
$ws.Range("T22").Value = 0.05857246964383642
$ws.Range("G23").Value = 67.39800266666667
$ws.Range("H23").Value = 202.194008
$ws.Range("I23").Value = 0.3905845264378918
$ws.Range("J23").Value = 0.3937305443979143
$ws.Range("O23").Value = 0.2991130341144489
$ws.Range("P23").Value = 0.3657579639239645
$ws.Range("Q23").Value = 99.56515972939113
$ws.Range("R23").Value = 896.08643756452
$ws.Range("S23").Value = 0.116828922780993
$ws.Range("T23").Value = 0.1440100822536552
$ws.Range("G24").Value = 67.39800266666667
$ws.Range("H24").Value = 202.194008
$ws.Range("I24").Value = 0.3905845264378918
$ws.Range("J24").Value = 0.3937305443979143
$ws.Range("M24").Value = 0.111967
$ws.Range("N24").Value = 0.335901
$ws.Range("O24").Value = 0.02267070427625646
$ws.Range("P24").Value = 0.02772193014375004
$ws.Range("Q24").Value = 7.546352164578667
$ws.Range("R24").Value = 67.91716948120801
$ws.Range("S24").Value = 0.008854826293755115
$ws.Range("T24").Value = 0.01091497064725965
$ws.Range("G25").Value = 67.39800266666667
$ws.Range("H25").Value = 202.194008
$ws.Range("I25").Value = 0.3905845264378918
$ws.Range("J25").Value = 0.3937305443979143
$ws.Range("M25").Value = 2.6997255
$ws.Range("N25").Value = 5.399451
$ws.Range("O25").Value = 0.5466314042313235
$ws.Range("P25").Value = 0.4456170223863617
$ws.Range("Q25").Value = 181.956106448268
$ws.Range("R25").Value = 1091.736638689608
$ws.Range("S25").Value = 0.2135057681577713
$ws.Range("T25").Value = 0.1754530328171598
$ws.Range("G26").Value = 67.39800266666667
$ws.Range("H26").Value = 202.194008
$ws.Range("I26").Value = 0.3905845264378918
$ws.Range("J26").Value = 0.3937305443979143
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.04903366666666667
$ws.Range("N26").Value = 0.147101
$ws.Range("O26").Value = 0.009928173091897913
$ws.Range("P26").Value = 0.01214025455737189
$ws.Range("Q26").Value = 3.304771196756445
$ws.Range("R26").Value = 29.742940770808
$ws.Range("S26").Value = 0.003877790785492366
$ws.Range("T26").Value = 0.004779989036003294
